$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value, even if it looks numeric,
# without leaving a residual number-format style on the cell (matches the
# original file where these data cells carry no explicit style).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "66.941.85"
Set-TextValue $ws.Range("E2") "  -1.75%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.604.11"
Set-TextValue $ws.Range("E3") "  -1.10%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.17%  "

# Row 5
Set-TextValue $ws.Range("D5") "586.48"
Set-TextValue $ws.Range("E5") "  -0.10%  "

# Row 6
Set-TextValue $ws.Range("D6") "182.06"
Set-TextValue $ws.Range("E6") "  +2.05%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.609"
Set-TextValue $ws.Range("E7") "  -2.67%  "

# Row 8
Set-TextValue $ws.Range("E8") "  +0.15%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.669"
Set-TextValue $ws.Range("E9") "  -5.60%  "

# Row 10
Set-TextValue $ws.Range("D10") "53.43"
Set-TextValue $ws.Range("E10") "  -2.94%  "

# Row 11
Set-TextValue $ws.Range("E11") "  -10.42%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0000251"
Set-TextValue $ws.Range("E12") "  -12.90%  "

# Row 13
Set-TextValue $ws.Range("D13") "9.88"
Set-TextValue $ws.Range("E13") "  -6.33%  "

# Row 14
Set-TextValue $ws.Range("D14") "4.184.18"
Set-TextValue $ws.Range("E14") "  -0.98%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.605.35"
Set-TextValue $ws.Range("E15") "  -1.10%  "

# Row 16
Set-TextValue $ws.Range("E16") "  -0.18%  "

# Row 17
Set-TextValue $ws.Range("D17") "66.777.40"
Set-TextValue $ws.Range("E17") "  -1.63%  "

# Row 18
Set-TextValue $ws.Range("D18") "18.27"
Set-TextValue $ws.Range("E18") "  -4.86%  "

# Row 19
Set-TextValue $ws.Range("D19") "12.13"
Set-TextValue $ws.Range("E19") "  -4.14%  "

# Row 20
Set-TextValue $ws.Range("D20") "1.05"
Set-TextValue $ws.Range("E20") "  -5.18%  "

# Row 21
Set-TextValue $ws.Range("D21") "390.57"
Set-TextValue $ws.Range("E21") "  -4.12%  "

# Row 22
Set-TextValue $ws.Range("D22") "4.29"
Set-TextValue $ws.Range("E22") "  -5.56%  "

# Row 23
Set-TextValue $ws.Range("D23") "84.54"
Set-TextValue $ws.Range("E23") "  -3.74%  "

# Row 24
Set-TextValue $ws.Range("D24") "2.86"
Set-TextValue $ws.Range("E24") "  -4.38%  "

# Row 25
Set-TextValue $ws.Range("D25") "12.18"
Set-TextValue $ws.Range("E25") "  -3.42%  "

# Row 26
Set-TextValue $ws.Range("D26") "6.03"
Set-TextValue $ws.Range("E26") "  -0.16%  "

# Row 27
Set-TextValue $ws.Range("D27") "10.22"
Set-TextValue $ws.Range("E27") "  -4.12%  "

# Row 28
Set-TextValue $ws.Range("D28") "3.60"
Set-TextValue $ws.Range("E28") "  -7.20%  "

# Row 29
Set-TextValue $ws.Range("D29") "8.91"
Set-TextValue $ws.Range("E29") "  -5.37%  "

# Row 30
Set-TextValue $ws.Range("D30") "31.02"
Set-TextValue $ws.Range("E30") "  -4.18%  "

# Row 31
Set-TextValue $ws.Range("D31") "6.72"
Set-TextValue $ws.Range("E31") "  -5.86%  "

# Row 32
Set-TextValue $ws.Range("D32") "11.86"
Set-TextValue $ws.Range("E32") "  -3.20%  "

# Row 33
Set-TextValue $ws.Range("D33") "64.84"
Set-TextValue $ws.Range("E33") "  +0.50%  "

# Row 34
Set-TextValue $ws.Range("D34") "595.45"
Set-TextValue $ws.Range("E34") "  -0.75%  "

# Row 35
Set-TextValue $ws.Range("E35") "  -3.52%  "

# Row 36
Set-TextValue $ws.Range("D36") "41.04"
Set-TextValue $ws.Range("E36") "  -3.45%  "

# Row 37
Set-TextValue $ws.Range("E37") "  +0.22%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.997"
Set-TextValue $ws.Range("E38") "  -0.35%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.371"
Set-TextValue $ws.Range("E39") "  -6.02%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0₃0734"
Set-TextValue $ws.Range("E40") "  -16.46%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.131"
Set-TextValue $ws.Range("E41") "  -3.71%  "

# Row 42
Set-TextValue $ws.Range("D42") "2.75"
Set-TextValue $ws.Range("E42") "  -8.17%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.0408"
Set-TextValue $ws.Range("E43") "  -6.07%  "

# Row 44
Set-TextValue $ws.Range("D44") "2.750.74"
Set-TextValue $ws.Range("E44") "  +2.16%  "

# Row 45
Set-TextValue $ws.Range("D45") "2.39"
Set-TextValue $ws.Range("E45") "  -10.47%  "

# Row 46
Set-TextValue $ws.Range("B46") "ApeXProtocol"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D46") "3.08"
Set-TextValue $ws.Range("E46") "  -0.61%  "

# Row 47
Set-TextValue $ws.Range("B47") "Stellar"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D47") "0.129"
Set-TextValue $ws.Range("E47") "  -3.57%  "

# Row 48
Set-TextValue $ws.Range("D48") "2.53"
Set-TextValue $ws.Range("E48") "  -6.13%  "

# Row 49
Set-TextValue $ws.Range("D49") "135.08"
Set-TextValue $ws.Range("E49") "  -3.49%  "

# Row 50
Set-TextValue $ws.Range("D50") "8.22"
Set-TextValue $ws.Range("E50") "  -7.88%  "

# Row 51
Set-TextValue $ws.Range("D51") "2.56"
Set-TextValue $ws.Range("E51") "  -6.26%  "
